# Apply updated odds/values per the FlashScore weekly games sheet diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
# Row 4
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 3.8
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("AE4").Value = 19
$ws.Range("AG4").Value = 12
$ws.Range("AM4").Value = 451
$ws.Range("AO4").Value = 8.5
$ws.Range("AP4").Value = 21
# Row 8
$ws.Range("N8").Value = 8.5
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.8
$ws.Range("X8").Value = 9.5
$ws.Range("AB8").Value = 34
$ws.Range("AC8").Value = 8.5
$ws.Range("AM8").Value = 351
$ws.Range("BA8").Value = 101
$ws.Range("BB8").Value = 251
# Row 9
$ws.Range("G9").Value = 1.5
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 6.5
$ws.Range("J9").Value = 2.1
$ws.Range("Z9").Value = 10
$ws.Range("AF9").Value = 67
$ws.Range("AS9").Value = 151
# Row 10
$ws.Range("G10").Value = 1.57
$ws.Range("H10").Value = 3.75
$ws.Range("I10").Value = 6.5
$ws.Range("J10").Value = 2.2
$ws.Range("L10").Value = 6
$ws.Range("O10").Value = 1.3
$ws.Range("P10").Value = 3.5
$ws.Range("Q10").Value = 2.03
$ws.Range("R10").Value = 1.87
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.75
$ws.Range("X10").Value = 7
$ws.Range("Z10").Value = 11
$ws.Range("AG10").Value = 15
$ws.Range("AH10").Value = 29
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 67
$ws.Range("AK10").Value = 51
$ws.Range("AL10").Value = 51
$ws.Range("AN10").Value = 3.4
$ws.Range("AO10").Value = 8
$ws.Range("AQ10").Value = 26
$ws.Range("AU10").Value = 9
$ws.Range("AV10").Value = 67
$ws.Range("AX10").Value = 34
$ws.Range("AY10").Value = 41
$ws.Range("AZ10").Value = 126
$ws.Range("BA10").Value = 151
# Row 11
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 7
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.75
$ws.Range("Q11").Value = 2.4
$ws.Range("R11").Value = 1.53
# Row 13
$ws.Range("O13").Value = 1.25
$ws.Range("P13").Value = 4
# Row 15
$ws.Range("G15").Value = 1.73
$ws.Range("H15").Value = 3.7
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 2.3
$ws.Range("L15").Value = 4.75
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11
$ws.Range("O15").Value = 1.25
$ws.Range("P15").Value = 4
$ws.Range("Q15").Value = 1.75
$ws.Range("R15").Value = 2.05
$ws.Range("S15").Value = 1.36
$ws.Range("T15").Value = 3
$ws.Range("U15").Value = 1.75
$ws.Range("V15").Value = 2
$ws.Range("W15").Value = 8
$ws.Range("X15").Value = 8.5
$ws.Range("AB15").Value = 23
$ws.Range("AC15").Value = 12
$ws.Range("AF15").Value = 41
$ws.Range("AG15").Value = 15
$ws.Range("AI15").Value = 15
$ws.Range("AK15").Value = 34
$ws.Range("AM15").Value = 201
$ws.Range("AN15").Value = 3.75
$ws.Range("AP15").Value = 19
$ws.Range("AS15").Value = 126
$ws.Range("AT15").Value = 3
$ws.Range("AU15").Value = 8
$ws.Range("AX15").Value = 23
$ws.Range("AY15").Value = 29
$ws.Range("BB15").Value = 201
# Row 20
$ws.Range("G20").Value = 2.2
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 3.1
$ws.Range("K20").Value = 1.8
$ws.Range("L20").Value = 4.75
$ws.Range("Z20").Value = 21
$ws.Range("AE20").Value = 21
$ws.Range("AG20").Value = 7.5
$ws.Range("AH20").Value = 17
$ws.Range("AN20").Value = 4
